$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Foglio1")

# Nuova riga (10): ricalcolo del "Dataset Misto" scalando le feature
# (Y invece di N), con la nuova accuracy ottenuta.
$ws.Range("A10").Value = "Dataset Misto"
$ws.Range("B10").Value = "0-23"
$ws.Range("C10").Value = 500
$ws.Range("D10").Value = "130 min"
$ws.Range("E10").Value = "Daniele(20 jobs)"
$ws.Range("F10").Value = 0.9998
$ws.Range("G10").Value = "Y"
$ws.Range("H10").Value = "Leap1 per 0-11, 21,22,23                                                 Leap2 per 12-20"

# Commento evidenziato con sottolineatura.
$ws.Range("I8").Font.Underline = $true
